# The edit reorders/swaps the content of the observation rows (9-18, minus
# the two rows 14 and 16 which are untouched) so that each row's taxon
# block (Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Ost, Nord) moves to a different row while the
# shared location/date/observer columns stay put (they are identical across
# the whole group already).
#
# Source row (content to copy) for each destination row:
#   9  <- 15
#   10 <- 11
#   11 <- 17
#   12 <- 18
#   13 <- 9
#   15 <- 12
#   17 <- 10
#   18 <- 13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$srcMap = @{
    9  = 15
    10 = 11
    11 = 17
    12 = 18
    13 = 9
    15 = 12
    17 = 10
    18 = 13
}

# Phase 1: snapshot the values of every row that participates in the
# permutation before any writes happen (it's a set of cycles, so reads must
# all happen before any write clobbers a still-needed source row).
$snapshot = @{}
foreach ($row in $srcMap.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Phase 2: write the snapshotted source-row values into each destination row.
foreach ($destRow in $srcMap.Keys) {
    $srcRow = $srcMap[$destRow]
    $rowValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowValues[$col]
    }
}
